# Update Name of Algo
# Applies updated numeric values (columns B and C) produced by a re-run of
# the KNN imputation algorithm for the terrestrial_mammals / combination_1_ABCD
# / BC / 20 / seed2 result dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.707
$ws.Range("B14").Value = 5.651999999999999
$ws.Range("C15").Value = -13.056
$ws.Range("B16").Value = 4.968999999999999
$ws.Range("B21").Value = 9.273999999999999
$ws.Range("C21").Value = -12.089
$ws.Range("C22").Value = -12.692
$ws.Range("B23").Value = 8.154
$ws.Range("C24").Value = -12.367
$ws.Range("B25").Value = 5.628
$ws.Range("B26").Value = 5.456
$ws.Range("C27").Value = -13.153
$ws.Range("C28").Value = -12.91
$ws.Range("B29").Value = 5.432
$ws.Range("C36").Value = -13.045
$ws.Range("C39").Value = -13.019
$ws.Range("B40").Value = 9.184999999999999
$ws.Range("C45").Value = -12.859
$ws.Range("C48").Value = -11.516
$ws.Range("C49").Value = -12.901
$ws.Range("C52").Value = -11.493
$ws.Range("B53").Value = 5.168
$ws.Range("C53").Value = -10.894
$ws.Range("C54").Value = -12.92
$ws.Range("B57").Value = 4.909
$ws.Range("C57").Value = -13.462
$ws.Range("B59").Value = 4.665000000000001
$ws.Range("B65").Value = 6.032999999999999
$ws.Range("B69").Value = 5.443
$ws.Range("C70").Value = -11.581
$ws.Range("C71").Value = -11.465
$ws.Range("B79").Value = 5.577
$ws.Range("B83").Value = 5.1
$ws.Range("C86").Value = -13.754
$ws.Range("C87").Value = -13.079
$ws.Range("C89").Value = -13.153
$ws.Range("B91").Value = 5.707
$ws.Range("B93").Value = 5.131000000000001
$ws.Range("B100").Value = 6.173999999999999
$ws.Range("C101").Value = -12.565
$ws.Range("B103").Value = 5.646
